$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text representation (values like "30.486.39"
# or "0.9993" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.486.39'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '1.911.66'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '239.56'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').Value = '0.9987'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').Value = '0.4755'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('D8').Value = '0.2845'
$ws.Range('E8').Value = '  -3.10%  '
$ws.Range('D9').Value = '0.06693'
$ws.Range('E9').Value = '  -5.39%  '
$ws.Range('E10').Value = '  -3.84%  '
$ws.Range('D11').Value = '101.24'
$ws.Range('E11').Value = '  -5.89%  '
$ws.Range('D12').Value = '1.916.78'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D13').Value = '0.07675'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').Value = '5.233'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '0.6707'
$ws.Range('E15').Value = '  -4.27%  '
$ws.Range('D16').Value = '30.501.81'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').Value = '256.78'
$ws.Range('E17').Value = '  -7.58%  '
$ws.Range('D18').Value = '0.9989'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '0.000007482'
$ws.Range('E19').Value = '  -4.17%  '
$ws.Range('D20').Value = '12.65'
$ws.Range('E20').Value = '  -4.83%  '
$ws.Range('D21').Value = '5.399'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').Value = '0.9994'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = '0.4519'
$ws.Range('E23').Value = '  -9.62%  '
$ws.Range('D24').Value = '6.310'
$ws.Range('E24').Value = '  -2.99%  '
$ws.Range('D25').Value = '168.64'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').Value = '9.351'
$ws.Range('E26').Value = '  -4.17%  '
$ws.Range('E27').Value = '  -3.50%  '
$ws.Range('D28').Value = '2.054'
$ws.Range('E28').Value = '  -5.30%  '
$ws.Range('D29').Value = '0.1008'
$ws.Range('E29').Value = '  -3.98%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '4.697'
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '1.369'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('D33').Value = '4.255'
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('D34').Value = '0.04715'
$ws.Range('E34').Value = '  -3.57%  '
$ws.Range('D35').Value = '0.7289'
$ws.Range('D36').Value = '1.112'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('D37').Value = '0.9981'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').Value = '2.710'
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').Value = '0.01919'
$ws.Range('E39').Value = '  -4.02%  '
$ws.Range('D40').Value = '2.597'
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('D41').Value = '6.235'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('D42').Value = '74.81'
$ws.Range('E42').Value = '  -4.78%  '
$ws.Range('D43').Value = '1.958'
$ws.Range('E43').Value = '  -7.55%  '
$ws.Range('D44').Value = '0.8616'
$ws.Range('E44').Value = '  -3.95%  '
$ws.Range('D45').Value = '105.41'
$ws.Range('E45').Value = '  -3.61%  '
$ws.Range('D46').Value = '0.4245'
$ws.Range('D47').Value = '0.9981'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').Value = '993.21'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = '7.407'
$ws.Range('E49').Value = '  -5.55%  '
$ws.Range('D50').Value = '0.1199'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('D51').Value = '34.87'
$ws.Range('E51').Value = '  -2.95%  '
